$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C to fit the longer description text
$ws.Columns.Item(3).ColumnWidth = 102.5

# New timesheet rows (14-19)
$ws.Range("A14").Value = 43072
$ws.Range("A14").NumberFormat = "d-mmm"
$ws.Range("B14").Value = 0.05555555555555555
$ws.Range("B14").NumberFormat = "h:mm"
$ws.Range("C14").Value = "Basic level design + destructables + Axe + Pick up (Key) + Cuttable Tree + Door locked (item requirement)"

$ws.Range("A15").Value = 43072
$ws.Range("A15").NumberFormat = "d-mmm"
$ws.Range("B15").Value = 0.052083333333333336
$ws.Range("B15").NumberFormat = "h:mm"
$ws.Range("C15").Value = "Created shop + shop panel"

$ws.Range("A16").Value = 43072
$ws.Range("A16").NumberFormat = "d-mmm"
$ws.Range("B16").Value = 0.020833333333333332
$ws.Range("B16").NumberFormat = "h:mm"
$ws.Range("C16").Value = "Show text above shopkeeper and locked door + Locked door remove key on unlock + breakable pot + enemy kill reward"

$ws.Range("A17").Value = 43072
$ws.Range("A17").NumberFormat = "d-mmm"
$ws.Range("B17").Value = 0.006944444444444444
$ws.Range("B17").NumberFormat = "h:mm"
$ws.Range("C17").Value = "Game over / complete screen"

$ws.Range("A18").Value = 43072
$ws.Range("A18").NumberFormat = "d-mmm"
$ws.Range("B18").Value = 0.017361111111111112
$ws.Range("B18").NumberFormat = "h:mm"
$ws.Range("C18").Value = "Update level design"

$ws.Range("A19").Value = 43072
$ws.Range("A19").NumberFormat = "d-mmm"
$ws.Range("B19").Value = 0.034722222222222224
$ws.Range("B19").NumberFormat = "h:mm"
$ws.Range("C19").Value = "Nav agent spawn bug fix + shopkeeper stock and bugfixes + Dev / God keys to heal and give items"

# Trailing blank (but time-formatted) row
$ws.Range("B20").NumberFormat = "h:mm"

# Move the selection the way it ended up after entering the new rows
$ws.Range("C21").Select() | Out-Null
